$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: RXNO_DEF
$ws.Range("F1").Value = "RXNO_DEF"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.LineStyle = 1

$ws.Range("F2").Value = '[''p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]'', locstr("Process, i.e., a physical entity with a temporal evolution that ''has a meaning for the ontologist''", ''en'')]'
$ws.Range("F3").Value = '[''B is a disposition means: b is a realizable entity and b’s bearer is some material entity and b is such that if it ceases to exist, then its bearer is physically changed, and b’s realization occurs when and because this bearer is in some special physical circumstances, and this realization occurs in virtue of the bearer’s physical make-up. [BFO]'']'
